$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.759.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -3.21%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.794.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.59%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.38%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'315.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.55%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.51%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.5375"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.59%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3812"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.64%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07418"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.57%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'41.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.53%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'1.086"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -3.02%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.29%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'6.195"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.18%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'7.450"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.73%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'20.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.49%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.789.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.35%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'88.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.58%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.00001058"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.71%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06496"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.80%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.28%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'17.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.11%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.928"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.03%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'27.797.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.04%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'11.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.89%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D26").Value = "'156.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.84%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'20.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.24%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.994.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.50%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.310"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.04%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'121.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.85%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.107"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.36%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.1096"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +5.10%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.648"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.35%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.503"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.20%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.06955"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +7.07%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.2196"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.96%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.02275"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.88%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'5.049"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.18%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'8.466"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -5.49%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'11.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.28%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.6101"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.59%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.417"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.53%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.162"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.38%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'13.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.01%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'3.680"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.24%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.5686"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.47%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'124.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.69%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "EOS"
$ws.Range("C48").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D48").Value = "'1.171"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.32%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.905"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.95%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.06780"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.65%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'71.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.87%  "
$ws.Range("E51").Style = "Normal"
